$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header / label text updates
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Reschedule Days 1"
$ws.Range("G2").Value = "start reschedule 1"
$ws.Range("H2").Value = "Duration 1"
$ws.Range("I2").Value = "Reschedule Days 2"
$ws.Range("J2").Value = "start reschedule 2"
$ws.Range("K2").Value = "Duration 2"

# Task name change (row 12 "User testing" -> "User testing, general")
$ws.Range("B12").Value = "User testing, general"

# ---------------------------------------------------------------------------
# 2. New formulas for columns I and K (rows 3-19), mirroring F/H formulas
# ---------------------------------------------------------------------------
$ws.Range("I3:I19").Formula = "=IF(ISBLANK(J3),,J3-G3-H3)"
$ws.Range("K3:K19").Formula = "=IF(I3<>0,D3,)"

# ---------------------------------------------------------------------------
# 3. Row 4 (Progress documentation) - End date pushed back, F4 becomes a
#    plain value instead of formula
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = 43248
$ws.Range("F4").Value = 33

# ---------------------------------------------------------------------------
# 4. Row 7 (User testing, presets) - H7 becomes plain value; new J7/K7
# ---------------------------------------------------------------------------
$ws.Range("H7").Value = 4
$ws.Range("J7").Value = 43248
$ws.Range("J7").NumberFormat = "m/d/yyyy"
$ws.Range("K7").Value = 4

# ---------------------------------------------------------------------------
# 5. Row 8 (Progressive Preset Change) - new J8/K8 (K8 formula)
# ---------------------------------------------------------------------------
$ws.Range("J8").Value = 43248
$ws.Range("J8").NumberFormat = "m/d/yyyy"
$ws.Range("K8").Formula = "=IF(I8<>0,D8,)"

# ---------------------------------------------------------------------------
# 6. Row 10 (Evaluate depth map) - new G10 date; F10/H10 recompute;
#    I10 own formula; J10 blank date-styled cell
# ---------------------------------------------------------------------------
$ws.Range("G10").Value = 43248
$ws.Range("G10").NumberFormat = "m/d/yyyy"
$ws.Range("I10").Formula = "=IF(ISBLANK(J10),,J10-G10-H10)"
$ws.Range("J10").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# 7. Row 11 (Neural Network) - H11 becomes plain value; new J11/K11
# ---------------------------------------------------------------------------
$ws.Range("H11").Value = 7
$ws.Range("J11").Value = 43248
$ws.Range("J11").NumberFormat = "m/d/yyyy"
$ws.Range("K11").Value = 14

# ---------------------------------------------------------------------------
# 8. Row 12 (User testing, general) - new J12/K12
# ---------------------------------------------------------------------------
$ws.Range("J12").Value = 43264
$ws.Range("J12").NumberFormat = "m/d/yyyy"
$ws.Range("K12").Value = 4

# ---------------------------------------------------------------------------
# 9. Row 13 (Mock up blur) - new G13 date; F13/H13 recompute; J13 blank
# ---------------------------------------------------------------------------
$ws.Range("G13").Value = 43262
$ws.Range("G13").NumberFormat = "m/d/yyyy"
$ws.Range("J13").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# 10. Row 14 (Design pipeline) - new G14 date; F14/H14 recompute
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = 43264
$ws.Range("G14").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# 11. Row 19 (Integration with the viewer) - new G19 date; F19/H19 recompute
# ---------------------------------------------------------------------------
$ws.Range("G19").Value = 43262
$ws.Range("G19").NumberFormat = "m/d/yyyy"

# ---------------------------------------------------------------------------
# 12. Row 31 (new) / Row 32 formula source change
# ---------------------------------------------------------------------------
$ws.Range("E31").Value = 43246
$ws.Range("E31").NumberFormat = "m/d/yyyy"
$ws.Range("E32").Formula = "=E31-F32"

# ---------------------------------------------------------------------------
# 13. Column widths for the new columns G, I, J
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 14.592447916666666
$ws.Columns.Item(9).ColumnWidth = 16.736979166666668
$ws.Columns.Item(10).ColumnWidth = 13.307291666666666

# ---------------------------------------------------------------------------
# 14. Selection / view bookkeeping
# ---------------------------------------------------------------------------
$ws.Range("H31").Select() | Out-Null

# ---------------------------------------------------------------------------
# 15. Move / resize the chart to its new anchor position
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 709.1173077632874
$co.Top = 337.4999212598425
$co.Width = 695.0946062992126
$co.Height = 352.2122047244095
